$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.204.70"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "1.872.47"
$ws.Range("E3").Value = "  +3.66%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'311.93"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("D7").Value = "'0.5020"
$ws.Range("E7").Value = "  -1.92%  "
$ws.Range("D8").Value = "'0.3933"
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("D9").Value = "'0.09858"
$ws.Range("E9").Value = "  +26.36%  "
$ws.Range("D10").Value = "'1.141"
$ws.Range("E10").Value = "  +3.16%  "
$ws.Range("D11").Value = "'41.25"
$ws.Range("E11").Value = "  +0.47%  "
$ws.Range("D12").Value = "'6.477"
$ws.Range("E12").Value = "  +1.81%  "
$ws.Range("D13").Value = "'21.00"
$ws.Range("E13").Value = "  +3.09%  "
$ws.Range("D14").Value = "1.868.80"
$ws.Range("E14").Value = "  +3.55%  "
$ws.Range("D15").Value = "'1.002"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("D16").Value = "'7.399"
$ws.Range("E16").Value = "  +1.15%  "
$ws.Range("D17").Value = "'0.00001134"
$ws.Range("E17").Value = "  +5.50%  "
$ws.Range("D18").Value = "'93.58"
$ws.Range("E18").Value = "  +1.24%  "
$ws.Range("D19").Value = "'0.06637"
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("D20").Value = "'17.43"
$ws.Range("E20").Value = "  +0.87%  "
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("D22").Value = "'6.123"
$ws.Range("E22").Value = "  +2.07%  "
$ws.Range("D23").Value = "28.284.71"
$ws.Range("E23").Value = "  +0.29%  "
$ws.Range("D24").Value = "'11.34"
$ws.Range("E24").Value = "  +1.95%  "
$ws.Range("D25").Value = "'2.268"
$ws.Range("E25").Value = "  +1.41%  "
$ws.Range("D26").Value = "'2.568"
$ws.Range("E26").Value = "  +4.42%  "
$ws.Range("D27").Value = "'21.31"
$ws.Range("E27").Value = "  +4.14%  "
$ws.Range("D28").Value = "2.082.76"
$ws.Range("E28").Value = "  +3.46%  "
$ws.Range("D29").Value = "'158.38"
$ws.Range("E29").Value = "  -1.45%  "
$ws.Range("D30").Value = "'127.93"
$ws.Range("E30").Value = "  +0.32%  "
$ws.Range("E31").Value = "  -2.73%  "
$ws.Range("D32").Value = "'1.063"
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("D33").Value = "'5.633"
$ws.Range("E33").Value = "  +1.26%  "
$ws.Range("D34").Value = "'3.619"
$ws.Range("E34").Value = "  -1.03%  "
$ws.Range("D35").Value = "'0.06820"
$ws.Range("E35").Value = "  -3.61%  "
$ws.Range("D36").Value = "'9.522"
$ws.Range("E36").Value = "  +4.20%  "
$ws.Range("D37").Value = "'0.02392"
$ws.Range("E37").Value = "  +2.02%  "
$ws.Range("D38").Value = "'0.2187"
$ws.Range("E38").Value = "  +0.73%  "
$ws.Range("D39").Value = "'11.51"
$ws.Range("E39").Value = "  -0.61%  "
$ws.Range("D40").Value = "'5.019"
$ws.Range("E40").Value = "  +0.29%  "
$ws.Range("D41").Value = "'0.6310"
$ws.Range("E41").Value = "  +2.38%  "
$ws.Range("D42").Value = "'1.173"
$ws.Range("E42").Value = "  +1.27%  "
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("D44").Value = "'13.60"
$ws.Range("E44").Value = "  +3.58%  "
$ws.Range("D45").Value = "'0.6026"
$ws.Range("E45").Value = "  +1.46%  "
$ws.Range("D46").Value = "'3.666"
$ws.Range("E46").Value = "  -1.58%  "
$ws.Range("D47").Value = "'1.269"
$ws.Range("E47").Value = "  -2.69%  "
$ws.Range("D48").Value = "'124.83"
$ws.Range("E48").Value = "  -0.24%  "
$ws.Range("D49").Value = "'1.993"
$ws.Range("E49").Value = "  +3.79%  "
$ws.Range("D50").Value = "'1.201"
$ws.Range("E50").Value = "  -0.57%  "
$ws.Range("D51").Value = "'1.123"
$ws.Range("E51").Value = "  +5.67%  "
